# Auto-generated edit script applying the cryptos.xlsx price/volume refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.093.91"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.37"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.15"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.294"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("E10").Value = "  -2.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.046.85"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.35"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.800.21"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.065.48"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.622"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.00"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.58"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0777"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.83"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.17"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.26"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.62"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.404.37"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.650"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.35"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.16"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.35"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.919"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.36"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.38%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0509"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.55%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.03"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0138"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -8.21%  "
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.86"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.947.15"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("E51").Value = "  +0.16%  "
